# Apply the "Adding lab 22 and 23" edit to the cardiac arrest data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Blood Pressure(mmHg) Control reading updated
$ws.Range("B3").Value = "124/79"

# Cardiac Output(mL/min) row updated
$ws.Range("B4").Value = 5346

# Ventilation(L/min) row updated
$ws.Range("B5").Value = 6.5
$ws.Range("C5").Value = 54.3

# Symp Activity(Hz) row updated
$ws.Range("C6").Value = 14
